# Updates the cryptos price/volume snapshot (Price column D, Volume(1h) column E)
# to the latest scraped values, matching the "Updated cryptos list" GitHub Action commit.
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel keeps
# storing them as text (matching the source data, which uses localized "d.ddd.dd"-style
# strings rather than real numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.129.88'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '3.504.67'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('E7').Value = '  -1.52%  '
$ws.Range('D8').Value = '3.500.25'
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  -2.62%  '
$ws.Range('D11').Value = "'7.26"
$ws.Range('E11').Value = '  +7.42%  '
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').Value = "'46.10"
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('E14').Value = '  -1.62%  '
$ws.Range('D15').Value = '4.074.97'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').Value = "'8.34"
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('D17').Value = "'613.53"
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range('D18').Value = '3.506.47'
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('D19').Value = '70.150.62'
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').Value = "'17.46"
$ws.Range('E22').Value = '  -1.39%  '
$ws.Range('D23').Value = "'9.14"
$ws.Range('E23').Value = '  -9.03%  '
$ws.Range('D24').Value = "'98.50"
$ws.Range('E24').Value = '  +1.42%  '
$ws.Range('D25').Value = "'15.56"
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('E26').Value = '  -3.77%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('D29').Value = "'33.83"
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('D30').Value = "'8.97"
$ws.Range('E30').Value = '  -2.95%  '
$ws.Range('E31').Value = '  -4.61%  '
$ws.Range('E32').Value = '  -5.18%  '
$ws.Range('E33').Value = '  -5.01%  '
$ws.Range('D34').Value = "'6.81"
$ws.Range('E34').Value = '  -3.48%  '
$ws.Range('D35').Value = "'629.69"
$ws.Range('E35').Value = '  +11.00%  '
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('D37').Value = "'10.76"
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').Value = "'0.0482"
$ws.Range('E38').Value = '  +6.04%  '
$ws.Range('E39').Value = '  -4.53%  '
$ws.Range('D40').Value = "'56.82"
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('D42').Value = "'0.144"
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D43').Value = '3.363.03'
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').Value = '0.0₃0734'
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('E45').Value = '  -5.98%  '
$ws.Range('E46').Value = '  -4.76%  '
$ws.Range('E47').Value = '  -3.77%  '
$ws.Range('E48').Value = '  -3.89%  '
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('D50').Value = "'132.75"
$ws.Range('E50').Value = '  -1.07%  '
